$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.162.53"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.576.88"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.80"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  -2.94%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.798.16"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.05"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.569.95"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.37"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "26.157.90"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "0.0₃0726"
$ws.Range("E18").Value = "  -1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.80"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.00"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").Value = "1.280.90"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.612"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.11"
$ws.Range("E39").Value = "  -9.82%  "
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.765"
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.25"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "1.711.61"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.65"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  -4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("E51").Value = "  -1.64%  "
